$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("AK_1.png",0.563,0.002,0.981,0.429,"Tidak Diketahui","Salah")
    ,@("AK_2.png",0.874,0.003,0.993,0.571,"Akhlak Kamiswara","Benar")
    ,@("AK_3.png",0.868,0.003,0.985,0.857,"Akhlak Kamiswara","Benar")
    ,@("AK_4.png",0.547,0.002,0.986,0.571,"Muhammad Iqbal Baqi","Salah")
    ,@("AK_5.png",0.536,0.002,0.992,0.571,"Akhlak Kamiswara","Benar")
    ,@("MIB_1.png",1.624,0.005,0.98,0.714,"Muhammad Iqbal Baqi","Benar")
    ,@("MIB_2.png",1.191,0.004,0.982,0.429,"Tidak Diketahui","Salah")
    ,@("MIB_3.png",1.759,0.006,1.0,1,"Muhammad Iqbal Baqi","Benar")
    ,@("MIB_4.png",1.175,0.004,0.989,0.857,"Muhammad Iqbal Baqi","Benar")
    ,@("MIB_5.png",1.578,0.005,0.99,0.429,"Tidak Diketahui","Salah")
    ,@("AAH_1.png",0.879,0.003,0.984,0.429,"Tidak Diketahui","Salah")
    ,@("AAH_2.png",1.112,0.004,0.995,0.857,"Andrea Ayunove Hutami","Benar")
    ,@("AAH_3.png",0.9,0.003,0.979,0.714,"Andrea Ayunove Hutami","Benar")
    ,@("TI_1.png",0.95,0.003,0.966,0.429,"Tidak Diketahui","Salah")
    ,@("TI_2.png",0.955,0.003,0.994,0.429,"Tidak Diketahui","Salah")
    ,@("TI_3.png",0.749,0.002,0.994,1,"Toni Ismail","Benar")
    ,@("TI_4.png",0.747,0.003,0.99,1,"Toni Ismail","Benar")
    ,@("TI_5.png",1.083,0.004,0.996,0.571,"Toni Ismail","Benar")
    ,@("RAS_1.png",0.684,0.002,0.984,0.429,"Tidak Diketahui","Salah")
    ,@("RAS_2.png",1.051,0.004,0.973,0.857,"Ridha Ayu Salsabila","Benar")
    ,@("RAS_3.png",0.576,0.002,0.987,0.571,"Muhammad Iqbal Baqi","Salah")
    ,@("RAS_4.png",1.246,0.004,0.979,0.429,"Tidak Diketahui","Salah")
    ,@("RAS_5.png",1.189,0.004,0.98,0.857,"Ridha Ayu Salsabila","Benar")
    ,@("RR_1.png",1.516,0.005,0.985,0.429,"Tidak Diketahui","Salah")
    ,@("RR_2.png",1.355,0.004,0.99,0.429,"Tidak Diketahui","Salah")
    ,@("RR_3.png",1.144,0.004,0.92,0.429,"Tidak Diketahui","Salah")
    ,@("RR_4.png",1.376,0.005,0.991,0.571,"Rafiqo Rapitasari","Benar")
    ,@("RR_5.png",1.575,0.005,0.985,1.0,"Rafiqo Rapitasari","Benar")
    ,@("AR_1.png",0.9,0.003,0.976,1,"Arizli Romadhon","Benar")
    ,@("GA_1.png",1.66,0.005,0.986,1.0,"Gege Ardiyansyah","Benar")
    ,@("GA_2.png",0.709,0.002,0.99,0.571,"Muhammad Iqbal Baqi","Salah")
    ,@("GA_3.png",0.749,0.002,0.995,0.429,"Tidak Diketahui","Salah")
    ,@("FY_1.png",1.176,0.004,0.979,0.429,"Tidak Diketahui","Salah")
    ,@("FY_2.png",1.689,0.006,0.985,0.429,"Tidak Diketahui","Salah")
    ,@("FY_3.png",1.49,0.005,0.992,0.429,"Tidak Diketahui","Salah")
    ,@("FY_4.png",1.267,0.004,0.99,0.429,"Tidak Diketahui","Salah")
    ,@("TO_1.png",0.85,0.003,0.978,0.286,"Tidak Diketahui","Salah")
    ,@("TO_2.png",1.068,0.004,0.988,0.857,"Tiara Oktavian","Benar")
    ,@("TO_3.png",0.865,0.003,0.992,0.714,"Tiara Oktavian","Benar")
    ,@("TO_4.png",7.419,0.024,0.946,0.286,"Tidak Diketahui","Salah")
    ,@("TO_5.png",7.163,0.024,0.944,0.286,"Tidak Diketahui","Salah")
    ,@("TD_1.png",3.823,0.013,0.922,0.429,"Tidak Diketahui","Benar")
    ,@("TD_2.png",5.015,0.017,0.917,0.286,"Tidak Diketahui","Benar")
    ,@("TD_3.png",1.1,0.004,0.972,0.429,"Tidak Diketahui","Benar")
    ,@("TD_4.png",1.0,0.003,0.991,0.286,"Tidak Diketahui","Benar")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 1
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
